$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '29.482.40'
$ws.Cells.Item(2,5).Value = '  -0.05%  '
$ws.Cells.Item(3,4).Value = '1.902.12'
$ws.Cells.Item(3,5).Value = '  -0.58%  '
$ws.Cells.Item(4,5).Value = '  +0.31%  '
$ws.Cells.Item(5,4).NumberFormat = '@'
$ws.Cells.Item(5,4).Value = '325.99'
$ws.Cells.Item(5,4).Style = 'Normal'
$ws.Cells.Item(5,5).Value = '  -2.13%  '
$ws.Cells.Item(6,4).NumberFormat = '@'
$ws.Cells.Item(6,4).Value = '1.004'
$ws.Cells.Item(6,4).Style = 'Normal'
$ws.Cells.Item(6,5).Value = '  +0.24%  '
$ws.Cells.Item(7,4).NumberFormat = '@'
$ws.Cells.Item(7,4).Value = '0.4789'
$ws.Cells.Item(7,4).Style = 'Normal'
$ws.Cells.Item(7,5).Value = '  +2.39%  '
$ws.Cells.Item(8,4).NumberFormat = '@'
$ws.Cells.Item(8,4).Value = '0.4060'
$ws.Cells.Item(8,4).Style = 'Normal'
$ws.Cells.Item(8,5).Value = '  -1.07%  '
$ws.Cells.Item(9,4).NumberFormat = '@'
$ws.Cells.Item(9,4).Value = '0.08074'
$ws.Cells.Item(9,4).Style = 'Normal'
$ws.Cells.Item(9,5).Value = '  +0.36%  '
$ws.Cells.Item(10,5).Value = '  -1.08%  '
$ws.Cells.Item(11,4).NumberFormat = '@'
$ws.Cells.Item(11,4).Value = '23.37'
$ws.Cells.Item(11,4).Style = 'Normal'
$ws.Cells.Item(11,5).Value = '  +4.36%  '
$ws.Cells.Item(12,4).Value = '1.934.56'
$ws.Cells.Item(12,5).Value = '  +0.84%  '
$ws.Cells.Item(13,4).NumberFormat = '@'
$ws.Cells.Item(13,4).Value = '5.959'
$ws.Cells.Item(13,4).Style = 'Normal'
$ws.Cells.Item(13,5).Value = '  -0.31%  '
$ws.Cells.Item(14,4).NumberFormat = '@'
$ws.Cells.Item(14,4).Value = '7.084'
$ws.Cells.Item(14,4).Style = 'Normal'
$ws.Cells.Item(14,5).Value = '  -1.50%  '
$ws.Cells.Item(15,4).NumberFormat = '@'
$ws.Cells.Item(15,4).Value = '90.23'
$ws.Cells.Item(15,4).Style = 'Normal'
$ws.Cells.Item(15,5).Value = '  +0.32%  '
$ws.Cells.Item(16,5).Value = '  +0.33%  '
$ws.Cells.Item(17,4).NumberFormat = '@'
$ws.Cells.Item(17,4).Value = '0.06716'
$ws.Cells.Item(17,4).Style = 'Normal'
$ws.Cells.Item(17,5).Value = '  +1.99%  '
$ws.Cells.Item(18,4).NumberFormat = '@'
$ws.Cells.Item(18,4).Value = '0.00001032'
$ws.Cells.Item(18,4).Style = 'Normal'
$ws.Cells.Item(18,5).Value = '  -0.22%  '
$ws.Cells.Item(19,4).NumberFormat = '@'
$ws.Cells.Item(19,4).Value = '17.61'
$ws.Cells.Item(19,4).Style = 'Normal'
$ws.Cells.Item(19,5).Value = '  -1.00%  '
$ws.Cells.Item(20,5).Value = '  +0.22%  '
$ws.Cells.Item(21,4).Value = '29.499.25'
$ws.Cells.Item(21,5).Value = '  +0.14%  '
$ws.Cells.Item(22,4).NumberFormat = '@'
$ws.Cells.Item(22,4).Value = '5.549'
$ws.Cells.Item(22,4).Style = 'Normal'
$ws.Cells.Item(22,5).Value = '  -0.46%  '
$ws.Cells.Item(23,4).NumberFormat = '@'
$ws.Cells.Item(23,4).Value = '11.78'
$ws.Cells.Item(23,4).Style = 'Normal'
$ws.Cells.Item(23,5).Value = '  +2.34%  '
$ws.Cells.Item(24,4).NumberFormat = '@'
$ws.Cells.Item(24,4).Value = '2.162'
$ws.Cells.Item(24,4).Style = 'Normal'
$ws.Cells.Item(24,5).Value = '  -2.07%  '
$ws.Cells.Item(25,4).Value = '2.157.85'
$ws.Cells.Item(25,5).Value = '  +0.56%  '
$ws.Cells.Item(26,4).NumberFormat = '@'
$ws.Cells.Item(26,4).Value = '153.91'
$ws.Cells.Item(26,4).Style = 'Normal'
$ws.Cells.Item(26,5).Value = '  -1.01%  '
$ws.Cells.Item(27,4).NumberFormat = '@'
$ws.Cells.Item(27,4).Value = '19.89'
$ws.Cells.Item(27,4).Style = 'Normal'
$ws.Cells.Item(27,5).Value = '  +0.01%  '
$ws.Cells.Item(28,4).NumberFormat = '@'
$ws.Cells.Item(28,4).Value = '6.103'
$ws.Cells.Item(28,4).Style = 'Normal'
$ws.Cells.Item(28,5).Value = '  +5.87%  '
$ws.Cells.Item(29,4).NumberFormat = '@'
$ws.Cells.Item(29,4).Value = '2.092'
$ws.Cells.Item(29,4).Style = 'Normal'
$ws.Cells.Item(29,5).Value = '  -2.41%  '
$ws.Cells.Item(30,4).NumberFormat = '@'
$ws.Cells.Item(30,4).Value = '118.45'
$ws.Cells.Item(30,4).Style = 'Normal'
$ws.Cells.Item(30,5).Value = '  +0.89%  '
$ws.Cells.Item(31,4).NumberFormat = '@'
$ws.Cells.Item(31,4).Value = '1.033'
$ws.Cells.Item(31,4).Style = 'Normal'
$ws.Cells.Item(31,5).Value = '  -3.47%  '
$ws.Cells.Item(32,4).NumberFormat = '@'
$ws.Cells.Item(32,4).Value = '0.09495'
$ws.Cells.Item(32,4).Style = 'Normal'
$ws.Cells.Item(32,5).Value = '  +0.35%  '
$ws.Cells.Item(33,4).NumberFormat = '@'
$ws.Cells.Item(33,4).Value = '5.482'
$ws.Cells.Item(33,4).Style = 'Normal'
$ws.Cells.Item(33,5).Value = '  +1.27%  '
$ws.Cells.Item(34,4).NumberFormat = '@'
$ws.Cells.Item(34,4).Value = '3.549'
$ws.Cells.Item(34,4).Style = 'Normal'
$ws.Cells.Item(34,5).Value = '  -0.46%  '
$ws.Cells.Item(35,4).NumberFormat = '@'
$ws.Cells.Item(35,4).Value = '1.391'
$ws.Cells.Item(35,4).Style = 'Normal'
$ws.Cells.Item(35,5).Value = '  -2.64%  '
$ws.Cells.Item(36,4).NumberFormat = '@'
$ws.Cells.Item(36,4).Value = '0.06084'
$ws.Cells.Item(36,4).Style = 'Normal'
$ws.Cells.Item(36,5).Value = '  -0.61%  '
$ws.Cells.Item(37,5).Value = '  -0.64%  '
$ws.Cells.Item(38,5).Value = '  -0.42%  '
$ws.Cells.Item(39,4).NumberFormat = '@'
$ws.Cells.Item(39,4).Value = '0.5885'
$ws.Cells.Item(39,4).Style = 'Normal'
$ws.Cells.Item(39,5).Value = '  -0.15%  '
$ws.Cells.Item(40,4).NumberFormat = '@'
$ws.Cells.Item(40,4).Value = '7.928'
$ws.Cells.Item(40,4).Style = 'Normal'
$ws.Cells.Item(40,5).Value = '  -5.77%  '
$ws.Cells.Item(41,5).Value = '  -0.01%  '
$ws.Cells.Item(42,4).NumberFormat = '@'
$ws.Cells.Item(42,4).Value = '10.23'
$ws.Cells.Item(42,4).Style = 'Normal'
$ws.Cells.Item(42,5).Value = '  +0.17%  '
$ws.Cells.Item(43,4).NumberFormat = '@'
$ws.Cells.Item(43,4).Value = '1.288'
$ws.Cells.Item(43,4).Style = 'Normal'
$ws.Cells.Item(43,5).Value = '  +1.64%  '
$ws.Cells.Item(44,2).Value = 'RenderToken'
$ws.Cells.Item(44,3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(44,4).NumberFormat = '@'
$ws.Cells.Item(44,4).Value = '2.391'
$ws.Cells.Item(44,4).Style = 'Normal'
$ws.Cells.Item(44,5).Value = '  +1.63%  '
$ws.Cells.Item(45,2).Value = 'Cronos'
$ws.Cells.Item(45,3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(45,4).NumberFormat = '@'
$ws.Cells.Item(45,4).Value = '0.07764'
$ws.Cells.Item(45,4).Style = 'Normal'
$ws.Cells.Item(45,5).Value = '  +3.46%  '
$ws.Cells.Item(46,4).NumberFormat = '@'
$ws.Cells.Item(46,4).Value = '12.31'
$ws.Cells.Item(46,4).Style = 'Normal'
$ws.Cells.Item(46,5).Value = '  +1.10%  '
$ws.Cells.Item(47,4).NumberFormat = '@'
$ws.Cells.Item(47,4).Value = '0.5535'
$ws.Cells.Item(47,4).Style = 'Normal'
$ws.Cells.Item(47,5).Value = '  -0.68%  '
$ws.Cells.Item(48,4).NumberFormat = '@'
$ws.Cells.Item(48,4).Value = '1.925'
$ws.Cells.Item(48,4).Style = 'Normal'
$ws.Cells.Item(48,5).Value = '  -0.36%  '
$ws.Cells.Item(49,4).NumberFormat = '@'
$ws.Cells.Item(49,4).Value = '114.12'
$ws.Cells.Item(49,4).Style = 'Normal'
$ws.Cells.Item(49,5).Value = '  +0.68%  '
$ws.Cells.Item(50,4).NumberFormat = '@'
$ws.Cells.Item(50,4).Value = '72.45'
$ws.Cells.Item(50,4).Style = 'Normal'
$ws.Cells.Item(50,5).Value = '  +1.13%  '
$ws.Cells.Item(51,4).NumberFormat = '@'
$ws.Cells.Item(51,4).Value = '0.2921'
$ws.Cells.Item(51,4).Style = 'Normal'
$ws.Cells.Item(51,5).Value = '  -1.91%  '
